# ---------------------------------------------------------------------------
# "Add new user stories" - Requirements and Specifications.docx
#
# Before: Use Cases section has US1..US6
#   US1 Generate and Send Personalized PDF Report
#   US2 Log Emotional Event
#   US3 View Functionally Equivalent Situations
#   US4 Optimize PDF Report Content
#   US5 Import Qualtrics Assessment Results
#   US6 Collect User Experience Data
#
# After: Use Cases section has US1..US7
#   US1 Complete Survey                              (NEW)
#   US2 Generate Personalized Feedback                (NEW)
#   US3 Generate and Send Personalized PDF Report      (was US1, unchanged body)
#   US4 Log Emotional Event                            (was US2, unchanged body)
#   US5 View Functionally Equivalent Situations        (was US3, unchanged body)
#   US6 Import Qualtrics Assessment Results             (was US5, unchanged body)
#   US7 Collect User Experience Data                   (was US6, tiny wording tweak)
#   (US4 "Optimize PDF Report Content" is removed entirely)
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

# Helper: insert a brand-new paragraph immediately after $para, containing
# $text (or an empty paragraph if $text is empty) with the surrounding
# Arial/22 formatting, and return the new Paragraph object.
function Add-ParaAfter($para, [string]$text) {
    $para.Range.InsertParagraphAfter()
    $newPara = $para.Next()
    if ($text -ne "") {
        $newPara.Range.Text = $text
    }
    return $newPara
}

# Helper: find the Paragraph object whose text equals $exactText (trimming
# the trailing paragraph mark) - used so later steps work off a solid
# anchor rather than fragile character offsets.
function Get-ParaByText($doc, [string]$exactText) {
    foreach ($p in $doc.Paragraphs) {
        $t = $p.Range.Text
        $t = $t.TrimEnd([char]13, [char]7)
        if ($t -eq $exactText) {
            return $p
        }
    }
    return $null
}

# ===========================================================================
# Step 1 - Retitle US1 and insert the new US1/US2 stories ahead of what
#          will become US3 (the original US1 content).
# ===========================================================================
$d.Content.Find.Execute(
    "User Story US1: Generate and Send Personalized PDF Report",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "User Story US1: Complete Survey", 2) | Out-Null

$titlePara = Get-ParaByText $d "User Story US1: Complete Survey"

$newStoryText = @(
    "As a user, I want to be able to complete the survey so I can have my temperament logged.",
    "Feature: Survey Completion",
    "Scenario: User takes the survey",
    "Given the user is logged in",
    "When they fully answer the survey questions and click submit",
    "Then the system logs the survey answers",
    "",
    "User Story US2: Generate Personalized Feedback",
    "As an Admin, I want to be able to have the system generate personalized feedback based upon the user’s survey so that the user can see their survey’s results.",
    "Feature: Generated Personalized Feedback",
    "Scenario: The user has finished their survey",
    "Given the user has completed their survey",
    "When they submit the survey",
    "Then a personalized report will be generated based upon the user’s answers.",
    "",
    "User Story US3: Generate and Send Personalized PDF Report"
)

$curPara = $titlePara
foreach ($t in $newStoryText) {
    $curPara = Add-ParaAfter $curPara $t
}

# ===========================================================================
# Step 2 - Renumber old US2 -> US4, old US3 -> US5.
# ===========================================================================
$d.Content.Find.Execute(
    "User Story US2: Log Emotional Event",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "User Story US4: Log Emotional Event", 2) | Out-Null

$d.Content.Find.Execute(
    "User Story US3: View Functionally Equivalent Situations",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "User Story US5: View Functionally Equivalent Situations", 2) | Out-Null

# ===========================================================================
# Step 3 - Delete the entire old US4 "Optimize PDF Report Content" story
#          (title paragraph through its trailing blank paragraph), leaving
#          the blank paragraph that precedes it untouched.
# ===========================================================================
$delStartPara = Get-ParaByText $d "User Story US4: Optimize PDF Report Content"
$delEndPara = Get-ParaByText $d "User Story US5: Import Qualtrics Assessment Results"

$delRange = $d.Range($delStartPara.Range.Start, $delEndPara.Range.Start)
$delRange.Delete()

# ===========================================================================
# Step 4 - Renumber old US5 -> US6, old US6 -> US7.
# ===========================================================================
$d.Content.Find.Execute(
    "User Story US5: Import Qualtrics Assessment Results",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "User Story US6: Import Qualtrics Assessment Results", 2) | Out-Null

$d.Content.Find.Execute(
    "User Story US6: Collect User Experience Data",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "User Story US7: Collect User Experience Data", 2) | Out-Null

# ===========================================================================
# Step 5 - Small wording tweak in the final (US7) scenario.
# ===========================================================================
$d.Content.Find.Execute(
    "Then I am able to see the user’s results to the assessment and can query the database.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Then I can see the user’s results to the assessment and can query the database.", 2) | Out-Null

Write-Output "edit complete"
